# Needle calibration: re-sort the recorded curvature samples in
# chronological order (ascending by the "time (s)" column).
#
# The header row (row 1) stays put; the data rows are re-ordered so that
# column A is ascending, carrying columns B:D along with each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:D8")
$sortKey   = $ws.Range("A1")

# 1 = xlAscending, 1 = xlYes (range has a header row that should not move)
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1)
